$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 22225642
$ws.Range("I100").Value = 33335630
$ws.Range("J100").Value = 5664.6665
$ws.Range("K100").Value = 33335630
$ws.Range("L100").Value = 5664.6665
$ws.Range("M100").Value = -33335089
$ws.Range("N100").Value = -6746.6665
$ws.Range("H112").Value = 1350.5
$ws.Range("J112").Value = 1364.6154
$ws.Range("L112").Value = 4093.8462
$ws.Range("N112").Value = -6309.8462
$ws.Range("H135").Value = 939.2941
$ws.Range("I135").Value = 920.5
$ws.Range("J135").Value = 984.4
$ws.Range("K135").Value = 8284.5
$ws.Range("L135").Value = 8859.6
$ws.Range("M135").Value = -5749.5
$ws.Range("N135").Value = -13929.6
$ws.Range("H137").Value = 796970.1
$ws.Range("I137").Value = 2168505
$ws.Range("J137").Value = 2923.5789
$ws.Range("K137").Value = 6505515
$ws.Range("L137").Value = 8770.736699999999
$ws.Range("M137").Value = -6502965
$ws.Range("N137").Value = -13870.7367
$ws.Range("H138").Value = 3198.2246
$ws.Range("I138").Value = 1800.8125
$ws.Range("J138").Value = 3875.7576
$ws.Range("K138").Value = 5402.4375
$ws.Range("L138").Value = 11627.2728
$ws.Range("M138").Value = -262.4375
$ws.Range("N138").Value = -21907.2728
$ws.Range("H140").Value = 68361.2
$ws.Range("J140").Value = 68361.2
$ws.Range("L140").Value = 68361.2
$ws.Range("N140").Value = -78721.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 332.93332
$ws.Range("I2").Value = 332.93332
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 332.93332
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -219.93332
$ws.Range("N2").ClearContents()
$ws.Range("H12").Value = 16999
$ws.Range("J12").Value = 16999
$ws.Range("L12").Value = 16999
$ws.Range("N12").Value = -17345
$ws.Range("H61").Value = 3471.8
$ws.Range("I61").Value = 3471.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3471.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3259.8
$ws.Range("N61").ClearContents()
$ws.Range("H116").Value = 332.93332
$ws.Range("I116").Value = 332.93332
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 332.93332
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1961.06668
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 3611.5588
$ws.Range("I122").Value = 3351.76
$ws.Range("K122").Value = 10055.28
$ws.Range("M122").Value = -7605.280000000001
$ws.Range("H136").Value = 3471.8
$ws.Range("I136").Value = 3471.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10415.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7865.400000000001
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 39498.5
$ws.Range("J137").Value = 39498.5
$ws.Range("L137").Value = 39498.5
$ws.Range("N137").Value = -49698.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 332.93332
$ws.Range("I3").Value = 332.93332
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 332.93332
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -218.93332
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 20668.666
$ws.Range("I15").Value = 9999
$ws.Range("K15").Value = 9999
$ws.Range("M15").Value = -9772
$ws.Range("H99").Value = 3950
$ws.Range("I99").Value = 1742.5
$ws.Range("K99").Value = 1742.5
$ws.Range("M99").Value = -244.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 223838.72
$ws.Range("I31").Value = 676546.9
$ws.Range("J31").Value = 3005.4634
$ws.Range("K31").Value = 676546.9
$ws.Range("L31").Value = 3005.4634
$ws.Range("M31").Value = -676251.9
$ws.Range("N31").Value = -3595.4634
$ws.Range("H34").Value = 223838.72
$ws.Range("I34").Value = 676546.9
$ws.Range("J34").Value = 3005.4634
$ws.Range("K34").Value = 676546.9
$ws.Range("L34").Value = 3005.4634
$ws.Range("M34").Value = -676344.9
$ws.Range("N34").Value = -3409.4634
$ws.Range("H127").Value = 41891.11
$ws.Range("J127").Value = 41891.11
$ws.Range("L127").Value = 41891.11
$ws.Range("N127").Value = -51811.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1786405.5
$ws.Range("I113").Value = 631.549
$ws.Range("J113").Value = 6579798.5
$ws.Range("K113").Value = 1894.647
$ws.Range("L113").Value = 19739395.5
$ws.Range("M113").Value = 275.3530000000001
$ws.Range("N113").Value = -19743735.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2915.3845
$ws.Range("I80").Value = 2790.9092
$ws.Range("J80").Value = 3600
$ws.Range("K80").Value = 2790.9092
$ws.Range("L80").Value = 3600
$ws.Range("M80").Value = -1792.9092
$ws.Range("N80").Value = -5596
$ws.Range("H83").Value = 2915.3845
$ws.Range("I83").Value = 2790.9092
$ws.Range("J83").Value = 3600
$ws.Range("K83").Value = 13954.546
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -8962.546
$ws.Range("N83").Value = -27984
$ws.Range("H122").Value = 5985
$ws.Range("I122").Value = 5001.5
$ws.Range("J122").Value = 6771.8
$ws.Range("K122").Value = 15004.5
$ws.Range("L122").Value = 20315.4
$ws.Range("M122").Value = -12554.5
$ws.Range("N122").Value = -25215.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1340
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("H122").Value = 7462.75
$ws.Range("I122").Value = 2402
$ws.Range("K122").Value = 7206
$ws.Range("M122").Value = -4756
$ws.Range("H136").Value = 3301.795
$ws.Range("I136").Value = 1458.5294
$ws.Range("J136").Value = 4726.136
$ws.Range("K136").Value = 4375.5882
$ws.Range("L136").Value = 14178.408
$ws.Range("M136").Value = -1825.5882
$ws.Range("N136").Value = -19278.408

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 11180.8
$ws.Range("I17").Value = 997
$ws.Range("J17").Value = 17970
$ws.Range("K17").Value = 997
$ws.Range("L17").Value = 17970
$ws.Range("M17").Value = -825
$ws.Range("N17").Value = -18314
$ws.Range("H47").Value = 40069
$ws.Range("J47").Value = 40069
$ws.Range("L47").Value = 40069
$ws.Range("N47").Value = -41213
$ws.Range("H96").Value = 54438370
$ws.Range("I96").Value = 84209070
$ws.Range("J96").Value = 3402875.5
$ws.Range("K96").Value = 84209070
$ws.Range("L96").Value = 3402875.5
$ws.Range("M96").Value = -84207697
$ws.Range("N96").Value = -3405621.5
$ws.Range("H122").Value = 3341.182
$ws.Range("I122").Value = 1950.0667
$ws.Range("K122").Value = 5850.2001
$ws.Range("M122").Value = -3400.2001
$ws.Range("H139").Value = 46532.5
$ws.Range("J139").Value = 46532.5
$ws.Range("L139").Value = 46532.5
$ws.Range("N139").Value = -56812.5
